# fix: replace metadata with correct file
#
# Applies the changes described by the reference diff:
#  1. Rename worksheet "charts" -> "visualizations"
#  2. Update the "Package info" sheet (DHIS2 version / Created / Identifier)
#  3. Update "dashboardItems" sheet: normalize the "Content/item type" column
#     (Chart -> specific visualization type, Map -> MAP) and blank out the
#     Content UID / Content name for the now-untyped Map rows
#  4. Update the "visualizations" (formerly "charts") sheet: blank out the
#     placeholder " " Description cells
#  5. Update the "programs" sheet: bump the "Last updated" date

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename sheet "charts" -> "visualizations"
# ---------------------------------------------------------------------
$wsVisualizations = $wb.Worksheets.Item("charts")
$wsVisualizations.Name = "visualizations"

# ---------------------------------------------------------------------
# 2. Package info sheet
# ---------------------------------------------------------------------
$wsPkg = $wb.Worksheets.Item("Package info")
$wsPkg.Range("B5").Value = "DHIS2.34.4-aff07fb"
$wsPkg.Range("B6").Value = "20210406T141800"
$wsPkg.Range("B7").Value = "AEFI_TRACKER_V1.1.2_DHIS2.34.4-aff07fb_20210406T141800"

# ---------------------------------------------------------------------
# 3. dashboardItems sheet
# ---------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("dashboardItems")

# Column B ("Content/item type") normalization: "Chart" becomes the actual
# visualization type, "Map" becomes "MAP".
$typeChanges = @{
    3  = "SINGLE_VALUE"
    4  = "SINGLE_VALUE"
    5  = "MAP"
    6  = "PIE"
    7  = "PIE"
    8  = "SINGLE_VALUE"
    9  = "SINGLE_VALUE"
    10 = "PIE"
    11 = "PIE"
    12 = "PIE"
    13 = "STACKED_COLUMN"
    14 = "PIE"
    15 = "PIE"
    16 = "PIE"
    17 = "PIE"
    18 = "PIE"
    19 = "PIE"
    20 = "PIE"
    21 = "PIE"
    22 = "PIE"
    23 = "PIE"
    24 = "PIE"
    25 = "STACKED_BAR"
    26 = "BAR"
    27 = "COLUMN"
    30 = "COLUMN"
    31 = "COLUMN"
    32 = "COLUMN"
    33 = "COLUMN"
    34 = "MAP"
    35 = "MAP"
    36 = "MAP"
    37 = "SINGLE_VALUE"
    38 = "SINGLE_VALUE"
    39 = "SINGLE_VALUE"
    40 = "COLUMN"
    41 = "COLUMN"
    42 = "COLUMN"
    43 = "STACKED_BAR"
    44 = "STACKED_BAR"
    45 = "STACKED_BAR"
    47 = "PIE"
    48 = "PIE"
    49 = "PIE"
    50 = "PIE"
}

foreach ($row in $typeChanges.Keys) {
    $wsItems.Cells.Item($row, 2).Value = $typeChanges[$row]
}

# The four "Map" rows also lose their Content UID (col A) / Content name
# (col C) values in the reference file.
$mapRows = @(5, 34, 35, 36)
foreach ($row in $mapRows) {
    $wsItems.Cells.Item($row, 1).Value = ""
    $wsItems.Cells.Item($row, 3).Value = ""
}

# ---------------------------------------------------------------------
# 4. visualizations sheet (formerly "charts")
# ---------------------------------------------------------------------
$blankDescriptionRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,30,31,32,36,37,38,39,40,41,42)
foreach ($row in $blankDescriptionRows) {
    $wsVisualizations.Cells.Item($row, 2).Value = ""
}

# ---------------------------------------------------------------------
# 5. programs sheet
# ---------------------------------------------------------------------
$wsPrograms = $wb.Worksheets.Item("programs")
# The target value looks like a date ("2021-04-06"), but must be stored as
# literal text (matching the existing cell type) rather than converted to a
# date serial number, so force it in as text via the classic leading
# apostrophe convention.
$wsPrograms.Range("C2").Value = "'2021-04-06"
